$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 470.7143
$ws.Range("I92").Value = 278.8
$ws.Range("K92").Value = 278.8
$ws.Range("M92").Value = 969.2
$ws.Range("H112").Value = 2782.8333
$ws.Range("I112").Value = 925.25
$ws.Range("K112").Value = 2775.75
$ws.Range("M112").Value = -1667.75
$ws.Range("H131").Value = 5224.9473
$ws.Range("J131").Value = 6854.231
$ws.Range("L131").Value = 20562.693
$ws.Range("N131").Value = -30642.693
$ws.Range("H138").Value = 2707522.5
$ws.Range("I138").Value = 5408012
$ws.Range("J138").Value = 7033.4326
$ws.Range("K138").Value = 16224036
$ws.Range("L138").Value = 21100.2978
$ws.Range("M138").Value = -16218896
$ws.Range("N138").Value = -31380.2978
$ws.Range("H141").Value = 16356.143
$ws.Range("I141").Value = 2623.25
$ws.Range("J141").Value = 34666.668
$ws.Range("K141").Value = 7869.75
$ws.Range("L141").Value = 104000.004
$ws.Range("M141").Value = -2689.75
$ws.Range("N141").Value = -114360.004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 766.6667
$ws.Range("I11").Value = 150
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 150
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = -10
$ws.Range("N11").Value = -2280
$ws.Range("H20").Value = 38747.75
$ws.Range("I20").Value = 51997.75
$ws.Range("J20").Value = 5622.75
$ws.Range("K20").Value = 51997.75
$ws.Range("L20").Value = 5622.75
$ws.Range("M20").Value = -51750.75
$ws.Range("N20").Value = -6116.75
$ws.Range("H37").Value = 5262.25
$ws.Range("I37").Value = 1299.5
$ws.Range("K37").Value = 1299.5
$ws.Range("M37").Value = -1162.5
$ws.Range("H94").Value = 1463
$ws.Range("I94").Value = 1245
$ws.Range("K94").Value = 1245
$ws.Range("M94").Value = -794
$ws.Range("H134").Value = 502144.44
$ws.Range("I134").Value = 542215.6
$ws.Range("J134").Value = 7933.3335
$ws.Range("K134").Value = 1626646.8
$ws.Range("L134").Value = 23800.0005
$ws.Range("M134").Value = -1624111.8
$ws.Range("N134").Value = -28870.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("H31").Value = 3481.4314
$ws.Range("I31").Value = 2146.889
$ws.Range("J31").Value = 4982.7915
$ws.Range("K31").Value = 2146.889
$ws.Range("L31").Value = 4982.7915
$ws.Range("M31").Value = -1851.889
$ws.Range("N31").Value = -5572.7915
$ws.Range("H34").Value = 3481.4314
$ws.Range("I34").Value = 2146.889
$ws.Range("J34").Value = 4982.7915
$ws.Range("K34").Value = 2146.889
$ws.Range("L34").Value = 4982.7915
$ws.Range("M34").Value = -1944.889
$ws.Range("N34").Value = -5386.7915
$ws.Range("H103").Value = 6065.7144
$ws.Range("I103").Value = 3816.6667
$ws.Range("J103").Value = 19560
$ws.Range("K103").Value = 3816.6667
$ws.Range("L103").Value = 19560
$ws.Range("M103").Value = -2644.6667
$ws.Range("N103").Value = -21904

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 324.83334
$ws.Range("J17").Value = 300
$ws.Range("L17").Value = 900
$ws.Range("N17").Value = -1238
$ws.Range("H34").Value = 823.04346
$ws.Range("I34").Value = 224.75
$ws.Range("J34").Value = 949
$ws.Range("K34").Value = 674.25
$ws.Range("L34").Value = 2847
$ws.Range("M34").Value = -590.25
$ws.Range("N34").Value = -3015
$ws.Range("H39").Value = 3782.6924
$ws.Range("I39").Value = 700
$ws.Range("J39").Value = 4343.1816
$ws.Range("K39").Value = 2100
$ws.Range("L39").Value = 13029.5448
$ws.Range("M39").Value = -1806
$ws.Range("N39").Value = -13617.5448
$ws.Range("H55").Value = 5099.5
$ws.Range("I55").Value = 6000
$ws.Range("J55").Value = 4999.4443
$ws.Range("K55").Value = 18000
$ws.Range("L55").Value = 14998.3329
$ws.Range("M55").Value = -17823
$ws.Range("N55").Value = -15352.3329
$ws.Range("H131").Value = 824.25
$ws.Range("I131").Value = 353.33334
$ws.Range("J131").Value = 870.82416
$ws.Range("K131").Value = 1060.00002
$ws.Range("L131").Value = 2612.47248
$ws.Range("M131").Value = 3979.99998
$ws.Range("N131").Value = -12692.47248
$ws.Range("H139").Value = 2296.4285
$ws.Range("I139").Value = 1805
$ws.Range("J139").Value = 2951.6667
$ws.Range("K139").Value = 5415
$ws.Range("L139").Value = 8855.000100000001
$ws.Range("M139").Value = -275
$ws.Range("N139").Value = -19135.0001
$ws.Range("H141").Value = 3570.5
$ws.Range("I141").Value = 3541.8667
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 10625.6001
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -5445.6001
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6364.0527
$ws.Range("I70").Value = 5344
$ws.Range("K70").Value = 5344
$ws.Range("M70").Value = -5074
$ws.Range("H73").Value = 6364.0527
$ws.Range("I73").Value = 5344
$ws.Range("K73").Value = 5344
$ws.Range("M73").Value = -4408

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2499.3333
$ws.Range("I16").Value = 2499.3333
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2499.3333
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2329.3333
$ws.Range("N16").ClearContents()
$ws.Range("H132").Value = 3543.0344
$ws.Range("I132").Value = 3069.92
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 9209.76
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -6679.76
$ws.Range("N132").Value = -24560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 52309.2
$ws.Range("J133").Value = 52309.2
$ws.Range("L133").Value = 52309.2
$ws.Range("N133").Value = -62429.2
